$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary values ---
$ws.Range("E11").Value2 = 75920
$ws.Range("C13").Value2 = 3
$ws.Range("F13").Value2 = 2

# --- Row 16 (worker 1): RAUMIR ALFONSO ARIAS MARTINEZ ---
$ws.Range("C16").Value2 = "73226517"
$ws.Range("D16").Value2 = "RAUMIR ALFONSO ARIAS MARTINEZ"
$ws.Range("E16").Value2 = "2503"
$ws.Range("F16").Value2 = 32266
$ws.Range("G16").Value2 = 1423500

# --- Row 17 (worker 2): ORLANDO MORENO DE VOZ ---
$ws.Range("C17").Value2 = "1047456614"
$ws.Range("D17").Value2 = "ORLANDO MORENO DE VOZ"
$ws.Range("E17").Value2 = "2503"
$ws.Range("F17").Value2 = 34164
$ws.Range("G17").Value2 = 1423500

# --- Row 18 (worker 3, new last row): ALVARO SEGUNDO MONTES BOHORQUEZ ---
# First bring over the closing-border formatting from the old last data row (23)
# so row 18 becomes the visually closed bottom of the table.
$ws.Range("B23:J23").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("C18").Value2 = "9113901"
$ws.Range("D18").Value2 = "ALVARO SEGUNDO MONTES BOHORQUEZ"
$ws.Range("E18").Value2 = "2508"
$ws.Range("F18").Value2 = 9490
$ws.Range("G18").Value2 = 1423500

# --- Remove the now-obsolete worker rows 19-23 (old rows 19,20,21,22,23) ---
# This shifts the trailing "firma" rows (old 28,29) up to become rows 23,24.
$ws.Rows("19:23").Delete()
